$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.03864900479865655, 0.5850438662025559, 0.7245914535781601, 0.8512293777696821, 0.8736531364880852, 19)
    3  = @(-0.05426294671939138, 0.5223775699184166, 0.5969843308118188, 0.7726476110179975, 0.7930847273798034, 18)
    4  = @(0.04121462200618248, 0.4282064007022935, 0.4114353831211015, 0.6414322903636061, 0.6598070052540843, 17)
    5  = @(0.1513232839577969, 0.3552284090837819, 0.2325151131945527, 0.4821982094476842, 0.4728538908111392, 16)
    6  = @(0.1388059945772336, 0.3804343702871964, 0.1997170828207708, 0.4468971725361113, 0.4397036079933307, 15)
    7  = @(0.1640011163335129, 0.3317361490081165, 0.2137299874398079, 0.4623094066096946, 0.4485592918475749, 14)
    8  = @(0.1805696058411114, 0.3326462541290968, 0.1927252227720784, 0.4390048095090513, 0.4164890877762218, 13)
    9  = @(0.217823495412195, 0.3128752448572103, 0.2168929369963322, 0.4657176580250443, 0.4299418927031617, 12)
    10 = @(0.1640088941767343, 0.2728317493205706, 0.1256684404680174, 0.3544974477595254, 0.3296156480036049, 11)
    11 = @(0.1230033177014676, 0.3182986834576886, 0.2493480347795076, 0.499347609165707, 0.5101396525723122, 10)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
